$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellXml = @(
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>55 x 29</w:t><w:br/><w:t xml:space="preserve">  2    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>11 x 23</w:t><w:br/><w:t xml:space="preserve">  2    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>91 x 28</w:t><w:br/><w:t xml:space="preserve">  2    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>99 x 14</w:t><w:br/><w:t xml:space="preserve">  1    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>22 x 68</w:t><w:br/><w:t xml:space="preserve">  6    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>76 x 58</w:t><w:br/><w:t xml:space="preserve">  5    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>76 x 90</w:t><w:br/><w:t xml:space="preserve">  9    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>75 x 82</w:t><w:br/><w:t xml:space="preserve">  8    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>24 x 17</w:t><w:br/><w:t xml:space="preserve">  1    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>97 x 65</w:t><w:br/><w:t xml:space="preserve">  6    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>49 x 57</w:t><w:br/><w:t xml:space="preserve">  5    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>41 x 60</w:t><w:br/><w:t xml:space="preserve">  6    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>21 x 13</w:t><w:br/><w:t xml:space="preserve">  1    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>46 x 63</w:t><w:br/><w:t xml:space="preserve">  6    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>',
    '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>11 x 48</w:t><w:br/><w:t xml:space="preserve">  4    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>'
)

$rows = 5
$cols = 3
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $idx = ($r - 1) * $cols + ($c - 1)
        $payload = $cellXml[$idx]
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $payload + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $cell = $t.Cell($r, $c)
        $null = $cell.Range.InsertXML($xml)
    }
}

Write-Host "Done updating lattice multiplication table."